$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2:I10").Value = "d"
$ws.Range("H13").Select()
